# Guild.xlsx edit script
# - Adds a new "Ref" field/row to the Property1 sheet (row 7), pushing the
#   existing "Desc" row down to row 8.
# - Makes Property1 the active sheet/tab instead of Record.
# - Adds data validation for the new Ref row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Property1"
$ws2 = $wb.Worksheets.Item(2)   # "Record"

# --- 1. Insert a new row above the old "Desc" row (row 7) -----------------
$ws1.Rows.Item(7).Insert()

# --- 2. Populate the new row 7 ("Ref") -------------------------------------
$ws1.Range("A7").Value = "Ref"
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
foreach ($col in $cols) {
    $ws1.Range($col + "7").Value = $false
}

# --- 3. Formatting for the new row (wrap text, like the Desc row) ---------
$ws1.Range("A7:AA7").WrapText = $true

# --- 4. Data validation for the new row ------------------------------------
$ws1.Range("A7").Validation.Add(9, 1, 1)
$ws1.Range("B7:I7").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws1.Range("J7:AA7").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- 5. Make Property1 the active sheet/tab --------------------------------
$ws1.Activate()
$ws1.Range("V15").Select()
